$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.728.62"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").Value = "1.545.26"
$ws.Range("E3").Value = "  -1.81%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.75"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.40"
$ws.Range("E8").Value = "  -3.84%  "
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("E10").Value = "  -1.29%  "
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D13").Value = "1.547.91"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "26.705.35"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.20"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.63"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("E23").Value = "  -5.37%  "
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.80"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.86"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0459"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "1.332.45"
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.90"
$ws.Range("E34").Value = "  -0.62%  "
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("E39").Value = "  +1.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.74"
$ws.Range("E40").Value = "  +5.79%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.994"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.52"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.679.85"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("B47").Value = "mCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.25"
$ws.Range("E47").Value = "  -3.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.84"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "0.0₇0965"
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0951"
$ws.Range("E51").Value = "  -0.02%  "
